# ---------------------------------------------------------------------------
# Update "上海-漫展信息.xlsx" (Shanghai comic-con / events info workbook)
# to match the regenerated gh-pages data snapshot.
#
# Sheets (tab order):
#   1 = 展览      (Exhibitions)
#   2 = 演出      (Performances)
#   3 = 本地生活  (Local life)
#   4 = 全部类型  (All types - combined listing)
#
# Changes:
#   * Column F ("想去人数" / want-to-go count) bumped on many rows across
#     all four sheets.
#   * On sheet "演出" row 16, column G flips from the text "已售罄"
#     (sold out) to the numeric price 80, and F16 bumps 174 -> 176.
#   * On sheet "演出" a brand-new concert row is inserted at row 32
#     ("上海·《若月亮没来》王宇宙LETO巡回演唱会"), pushing the previous
#     rows 32-42 down to 33-43 (dimension A1:I42 -> A1:I43).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions) - column F updates only
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$sheet1F = @{
    4  = 580
    5  = 2552
    7  = 170
    9  = 244
    10 = 5207
    12 = 1447
    13 = 1373
    15 = 6923
    16 = 386
    17 = 45
    18 = 6
    20 = 4654
    22 = 67
    23 = 2294
    24 = 1247
    25 = 434
    27 = 212
    29 = 70
    30 = 150
    32 = 1264
    33 = 1979
    34 = 222
    35 = 510
    36 = 197
    37 = 1361
    39 = 84
    41 = 156
    42 = 1106
    43 = 2393
    45 = 60
    49 = 10
}
foreach ($row in $sheet1F.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1F[$row]
}

# ---------------------------------------------------------------------------
# Sheet 2: 演出 (Performances)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

# Plain column-F bumps on rows that sit above the new row 32 insert point,
# so they are unaffected by the later row shift.
$sheet2F = @{
    12 = 378
    13 = 264
    15 = 34
    20 = 132
    21 = 26
    28 = 274
}
foreach ($row in $sheet2F.Keys) {
    $ws2.Cells.Item($row, 6).Value = $sheet2F[$row]
}

# Row 16: F bumps 174 -> 176, and G flips from the "已售罄" text to the
# numeric price 80.
$ws2.Cells.Item(16, 6).Value = 176
$ws2.Cells.Item(16, 7).Value = 80

# Insert the new concert row at position 32; everything from the old
# row 32 ("吕思清小提琴独奏音乐会") through row 42 shifts down to 33-43.
$ws2.Rows.Item(32).Insert()

# Row 32's "A" column (sequence number, bold/centered/bordered style) lost
# its formatting on insert - clone it from the row below (which carries the
# original style forward) before writing the new index value.
$ws2.Cells.Item(33, 1).Copy()
$ws2.Cells.Item(32, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Cells.Item(32, 1).Value = 31
# Force column B to stay plain text ("2024-09-20") instead of Excel's
# automatic date-literal conversion - every other date cell in this sheet
# is stored as literal text, not a real date serial.
$ws2.Cells.Item(32, 2).NumberFormat = "@"
$ws2.Cells.Item(32, 2).Value = "2024-09-20"
$ws2.Cells.Item(32, 3).Value = "上海·《若月亮没来》王宇宙LETO巡回演唱会"
$ws2.Cells.Item(32, 4).Value = "嘉兴路街道瑞虹路188号瑞虹天地月亮湾3层 Modern Sky LAB摩登天空(瑞虹天地店)"
$ws2.Cells.Item(32, 5).Value = "2024.09.20 20:00-09.20 21:30"
$ws2.Cells.Item(32, 6).Value = 0
$ws2.Cells.Item(32, 7).Value = 168
$ws2.Cells.Item(32, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90197"
$ws2.Cells.Item(32, 9).Value = "//i0.hdslb.com/bfs/openplatform/202408/Bxe7VVNE1722499960950.jpeg"

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life) - column F updates only
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$sheet3F = @{
    6  = 1648
    8  = 1270
    10 = 1731
    11 = 2154
    12 = 582
    13 = 489
}
foreach ($row in $sheet3F.Keys) {
    $ws3.Cells.Item($row, 6).Value = $sheet3F[$row]
}

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (All types, combined listing) - column F updates only
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$sheet4F = @{
    3  = 1648
    5  = 580
    7  = 2552
    8  = 170
    9  = 1270
    10 = 2154
    11 = 5207
    12 = 582
    17 = 1447
    18 = 1373
    20 = 6923
    21 = 386
    22 = 489
    23 = 45
    24 = 4654
    25 = 2294
    26 = 1247
    27 = 434
    29 = 212
    30 = 70
    31 = 264
    33 = 150
    35 = 1979
    36 = 222
    37 = 510
    38 = 26
    39 = 1361
    42 = 156
    44 = 1106
    45 = 2393
    46 = 60
}
foreach ($row in $sheet4F.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4F[$row]
}

Write-Output "edit complete"
